$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '43.558.00'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.380.88'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +3.24%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '310.01'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '104.47'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +3.27%  '
$ws.Range("E7").Value = '  -4.90%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '53.43'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '2.752.23'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '15.60'
$ws.Range("D16").Style = 'Normal'
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '2.387.50'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +3.22%  '
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '43.549.18'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +1.00%  '
$ws.Range("E20").Value = '  +3.52%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '11.91'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -4.96%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.0₃0916'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '68.39'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '240.66'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +1.94%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.61'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '25.82'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +4.60%  '
$ws.Range("E29").Value = '  -3.09%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '36.57'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '9.52'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '160.92'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -3.48%  '
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '18.33'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +3.48%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '4.75'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +10.40%  '
$ws.Range("E38").Value = '  +5.61%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '3.11'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.0738'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  +5.69%  '
$ws.Range("E42").Value = '  -1.45%  '
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '2.60'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +13.72%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '2.035.60'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +2.46%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '19.76'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +3.41%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.0291'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '3.13'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +3.76%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '10.58'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +7.60%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '57.96'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +3.95%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '2.97'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.17%  '
